$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 100.8373843333333
$ws.Range("H2").Value = 302.512153
$ws.Range("I2").Value = 0.6551985585448407
$ws.Range("J2").Value = 0.6551985585448408
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.425703666666666
$ws.Range("N2").Value = 4.277111
$ws.Range("O2").Value = 0.04715501820393346
$ws.Range("P2").Value = 0.04715501820393346
$ws.Range("Q2").Value = 143.7642285811092
$ws.Range("R2").Value = 1293.878057229983
$ws.Range("S2").Value = 0.03089589995537292
$ws.Range("T2").Value = 0.03089589995537293

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 100.8373843333333
$ws.Range("H3").Value = 302.512153
$ws.Range("I3").Value = 0.6551985585448407
$ws.Range("J3").Value = 0.6551985585448408
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.61433933333333
$ws.Range("N3").Value = 61.843018
$ws.Range("O3").Value = 0.6818173855147049
$ws.Range("P3").Value = 0.6818173855147048
$ws.Range("Q3").Value = 2078.696058133084
$ws.Range("R3").Value = 18708.26452319775
$ws.Range("S3").Value = 0.4467257681800466
$ws.Range("T3").Value = 0.4467257681800466

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 100.8373843333333
$ws.Range("H4").Value = 302.512153
$ws.Range("I4").Value = 0.6551985585448407
$ws.Range("J4").Value = 0.6551985585448408
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.194356666666666
$ws.Range("N4").Value = 24.58307
$ws.Range("O4").Value = 0.2710275962813615
$ws.Range("P4").Value = 0.2710275962813615
$ws.Range("Q4").Value = 826.2974925610788
$ws.Range("R4").Value = 7436.67743304971
$ws.Range("S4").Value = 0.1775768904094211
$ws.Range("T4").Value = 0.1775768904094211

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 23.90796933333333
$ws.Range("H5").Value = 71.72390799999999
$ws.Range("I5").Value = 0.1553438454249564
$ws.Range("J5").Value = 0.1553438454249564
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.425703666666666
$ws.Range("N5").Value = 4.277111
$ws.Range("O5").Value = 0.04715501820393346
$ws.Range("P5").Value = 0.04715501820393346
$ws.Range("Q5").Value = 34.08567954108754
$ws.Range("R5").Value = 306.7711158697879
$ws.Range("S5").Value = 0.007325241858882845
$ws.Range("T5").Value = 0.007325241858882845

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 23.90796933333333
$ws.Range("H6").Value = 71.72390799999999
$ws.Range("I6").Value = 0.1553438454249564
$ws.Range("J6").Value = 0.1553438454249564
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 20.61433933333333
$ws.Range("N6").Value = 61.843018
$ws.Range("O6").Value = 0.6818173855147049
$ws.Range("P6").Value = 0.6818173855147048
$ws.Range("Q6").Value = 492.8469926082604
$ws.Range("R6").Value = 4435.622933474344
$ws.Range("S6").Value = 0.1059161345434442
$ws.Range("T6").Value = 0.1059161345434442

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 23.90796933333333
$ws.Range("H7").Value = 71.72390799999999
$ws.Range("I7").Value = 0.1553438454249564
$ws.Range("J7").Value = 0.1553438454249564
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.194356666666666
$ws.Range("N7").Value = 24.58307
$ws.Range("O7").Value = 0.2710275962813615
$ws.Range("P7").Value = 0.2710275962813615
$ws.Range("Q7").Value = 195.9104278930622
$ws.Range("R7").Value = 1763.19385103756
$ws.Range("S7").Value = 0.04210246902262932
$ws.Range("T7").Value = 0.04210246902262932

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.15819666666667
$ws.Range("H8").Value = 87.47459000000001
$ws.Range("I8").Value = 0.1894575960302029
$ws.Range("J8").Value = 0.1894575960302029
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.425703666666666
$ws.Range("N8").Value = 4.277111
$ws.Range("O8").Value = 0.04715501820393346
$ws.Range("P8").Value = 0.04715501820393346
$ws.Range("Q8").Value = 41.57094790105444
$ws.Range("R8").Value = 374.13853110949
$ws.Range("S8").Value = 0.008933876389677691
$ws.Range("T8").Value = 0.008933876389677691

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.15819666666667
$ws.Range("H9").Value = 87.47459000000001
$ws.Range("I9").Value = 0.1894575960302029
$ws.Range("J9").Value = 0.1894575960302029
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 20.61433933333333
$ws.Range("N9").Value = 61.843018
$ws.Range("O9").Value = 0.6818173855147049
$ws.Range("P9").Value = 0.6818173855147048
$ws.Range("Q9").Value = 601.0769604347356
$ws.Range("R9").Value = 5409.69264391262
$ws.Range("S9").Value = 0.1291754827912141
$ws.Range("T9").Value = 0.1291754827912141

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 29.15819666666667
$ws.Range("H10").Value = 87.47459000000001
$ws.Range("I10").Value = 0.1894575960302029
$ws.Range("J10").Value = 0.1894575960302029
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.194356666666666
$ws.Range("N10").Value = 24.58307
$ws.Range("O10").Value = 0.2710275962813615
$ws.Range("P10").Value = 0.2710275962813615
$ws.Range("Q10").Value = 238.9326632434778
$ws.Range("R10").Value = 2150.3939691913
$ws.Range("S10").Value = 0.05134823684931112
$ws.Range("T10").Value = 0.05134823684931112

